$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) / 1h-change (E) figures scraped for the crypto table.
# Helper: write a value as literal text so Excel does not silently
# reinterpret decimal-looking strings (e.g. '1.00', '6.87') as numbers,
# while leaving the cell's style untouched afterwards.
function Set-CellText($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $looksNumeric = $text -match '^\s*[+-]?[0-9]+(\.[0-9]+)?\s*$'
    if ($looksNumeric) {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

$updates = [ordered]@{
    'D2' = '43.069.68'
    'E2' = '  -0.15%  '
    'D3' = '2.299.21'
    'E3' = '  -0.39%  '
    'E4' = '  +0.02%  '
    'D5' = '300.58'
    'E5' = '  -0.35%  '
    'D6' = '98.31'
    'E6' = '  -1.82%  '
    'E7' = '  +1.29%  '
    'E8' = '  -0.02%  '
    'E9' = '  +0.52%  '
    'D10' = '36.19'
    'E10' = '  -0.69%  '
    'D11' = '0.0790'
    'E11' = '  -0.39%  '
    'E12' = '  +0.79%  '
    'E13' = '  -0.21%  '
    'D14' = '6.87'
    'D15' = '2.656.60'
    'E15' = '  -0.37%  '
    'D16' = '2.310.19'
    'E16' = '  +1.84%  '
    'D17' = '0.788'
    'E17' = '  -1.41%  '
    'D18' = '42.942.83'
    'E18' = '  -0.21%  '
    'D19' = '12.82'
    'E19' = '  -0.18%  '
    'E20' = '  +0.79%  '
    'D21' = '6.12'
    'E21' = '  -0.12%  '
    'D22' = '68.96'
    'E22' = '  +1.42%  '
    'D23' = '237.29'
    'E23' = '  +0.62%  '
    'E24' = '  -2.83%  '
    'E25' = '  -0.02%  '
    'E26' = '  -0.85%  '
    'D27' = '24.94'
    'E27' = '  -0.77%  '
    'D28' = '164.70'
    'E28' = '  -3.07%  '
    'E29' = '  -0.37%  '
    'E30' = '  -0.36%  '
    'D31' = '33.05'
    'E31' = '  -4.14%  '
    'D32' = '1.00'
    'E32' = '  +0.03%  '
    'E33' = '  +0.22%  '
    'D34' = '4.77'
    'E34' = '  +2.25%  '
    'D35' = '17.92'
    'E35' = '  +1.00%  '
    'D36' = '2.40'
    'E36' = '  -0.25%  '
    'E37' = '  +0.71%  '
    'E38' = '  -0.01%  '
    'E39' = '  -0.78%  '
    'D40' = '2.79'
    'E40' = '  -1.06%  '
    'E41' = '  +0.54%  '
    'D42' = '2.012.63'
    'E42' = '  +1.36%  '
    'E43' = '  -2.09%  '
    'D44' = '2.23'
    'E44' = '  -1.25%  '
    'E45' = '  +1.20%  '
    'D46' = '17.47'
    'E46' = '  -1.37%  '
    'E47' = '  -2.48%  '
    'D48' = '54.07'
    'E48' = '  -2.58%  '
    'D49' = '2.523.44'
    'E49' = '  -0.22%  '
    'D50' = '1.53'
    'E50' = '  -1.47%  '
    'D51' = '73.31'
    'E51' = '  +3.56%  '
}

foreach ($cellRef in $updates.Keys) {
    Set-CellText $cellRef $updates[$cellRef]
}
